# Update cryptos list data (prices and 1h volume change %) and reorder
# two coin pairs (Toncoin/ImmutableX and dogwifhat/WEMIXToken/VeChain/
# Stellar/THORChain) to reflect the refreshed ranking snapshot.
#
# Some Price values look like plain numbers (e.g. "207.18") but must stay
# as text cells exactly as authored (matching the original inline-string
# cells), so for those we force a text NumberFormat before assigning the
# value and then restore the "Normal" style so no stray per-cell style
# survives in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.847.38"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.563.92"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "566.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").Value = "3.562.64"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.03%  "
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "4.128.28"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "3.569.03"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "67.678.48"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "400.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.29%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "668.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.407"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.94%  "
$ws.Range("D40").Value = "0.0₃0743"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "3.153.76"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +17.26%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.46%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0408"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.75%  "
